$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts content/description/type columns right)
$ws.Columns("C:C").Insert()

# New header for the inserted "group" column
$ws.Range("C1").Value = "group"

# Update row 2: name -> "Gamma Ray", group -> "Clay Volume"
$ws.Range("B2").Value = "Gamma Ray"
$ws.Range("C2").Value = "Clay Volume"

# Update row 3: name -> "Density", group -> "Porosity"
$ws.Range("B3").Value = "Density"
$ws.Range("C3").Value = "Porosity"

# Column width for the new "group" column (maps to stored OOXML width of 10)
$ws.Columns("C:C").ColumnWidth = 9.140625
